$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-set NumberFormat to Text for price cells whose new values would
# otherwise be auto-converted to numbers by Excel's type inference.
$textCells = 'D5','D8','D11','D13','D16','D18','D20','D21','D25','D26','D28','D34','D39','D40','D42','D44','D45','D46','D48'
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '35.140.37'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.855.73'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '238.13'
$ws.Range("E5").Value = '  +3.25%  '
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '42.11'
$ws.Range("E8").Value = '  +5.13%  '
$ws.Range("E9").Value = '  +2.42%  '
$ws.Range("E10").Value = '  +1.25%  '
$ws.Range("D11").Value = '0.0989'
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("D12").Value = '2.123.22'
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '11.44'
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.859.00'
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").Value = '4.72'
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("D17").Value = '35.086.14'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = '69.91'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '240.83'
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").Value = '12.22'
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = '167.96'
$ws.Range("E25").Value = '  -3.43%  '
$ws.Range("D26").Value = '1.84'
$ws.Range("E26").Value = '  +22.88%  '
$ws.Range("E27").Value = '  +2.02%  '
$ws.Range("D28").Value = '17.65'
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("E33").Value = '  +23.29%  '
$ws.Range("D34").Value = '4.00'
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("E35").Value = '  +20.66%  '
$ws.Range("E36").Value = '  +10.92%  '
$ws.Range("E37").Value = '  +6.82%  '
$ws.Range("E38").Value = '  +6.65%  '
$ws.Range("D39").Value = '90.36'
$ws.Range("E39").Value = '  -2.54%  '
$ws.Range("D40").Value = '0.0201'
$ws.Range("E40").Value = '  +3.94%  '
$ws.Range("D41").Value = '1.341.87'
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").Value = '14.84'
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("E43").Value = '  +2.70%  '
$ws.Range("B44").Value = 'Gas'
$ws.Range("C44").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D44").Value = '12.54'
$ws.Range("E44").Value = '  +45.89%  '
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").Value = '2.41'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").Value = '0.0556'
$ws.Range("E46").Value = '  +6.45%  '
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("D48").Value = '6.50'
$ws.Range("E48").Value = '  +3.87%  '
$ws.Range("D49").Value = '2.038.23'
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("E51").Value = '  +0.28%  '
